$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a text value that looks like a date/time (e.g. "2026-01-31"
# or "08:44") into a cell while:
#  - preventing Excel from auto-converting it into a date/time serial number
#  - keeping the cell on the default ("Normal") style, matching the source
#    workbook which never sets an explicit style/number format on these cells
function Set-DateLikeText($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Column B (Taxonsorteringsordning) updates - simple numeric changes
# ---------------------------------------------------------------------
$ws.Range("B2").Value  = 79243
$ws.Range("B3").Value  = 79243
$ws.Range("B4").Value  = 79243
$ws.Range("B5").Value  = 79243
$ws.Range("B6").Value  = 79243
$ws.Range("B7").Value  = 79243
$ws.Range("B9").Value  = 79243
$ws.Range("B10").Value = 57884
$ws.Range("B11").Value = 57884
$ws.Range("B14").Value = 57884
$ws.Range("B15").Value = 79243
$ws.Range("B16").Value = 79243
$ws.Range("B17").Value = 79243
$ws.Range("B18").Value = 57884
$ws.Range("B19").Value = 79243
$ws.Range("B20").Value = 79243

# ---------------------------------------------------------------------
# Rows 12 and 13 are swapped (their underlying observation records traded
# places), and the B column also picks up the same 79239 -> 79243 change
# as the rest of the sheet.
# ---------------------------------------------------------------------

# Capture the current ("before") values of row 12 and row 13 first, since
# row 12 will be overwritten before row 13 is written.
$row12_A  = $ws.Range("A12").Value2
$row12_Q  = $ws.Range("Q12").Value2
$row12_R  = $ws.Range("R12").Value2
$row12_S  = $ws.Range("S12").Value2
$row12_AW = $ws.Range("AW12").Text
$row12_AX = $ws.Range("AX12").Text

$row13_A  = $ws.Range("A13").Value2
$row13_Q  = $ws.Range("Q13").Value2
$row13_R  = $ws.Range("R13").Value2
$row13_S  = $ws.Range("S13").Value2
$row13_Z  = $ws.Range("Z13").Text
$row13_AB = $ws.Range("AB13").Text
$row13_AW = $ws.Range("AW13").Text
$row13_AX = $ws.Range("AX13").Text

# --- New row 12 (takes former row 13 content) ---
$ws.Range("A12").Value = $row13_A
$ws.Range("B12").Value = 79243
$ws.Range("Q12").Value = $row13_Q
$ws.Range("R12").Value = $row13_R
$ws.Range("S12").Value = $row13_S

Set-DateLikeText $ws.Range("Z12") $row13_Z
Set-DateLikeText $ws.Range("AB12") $row13_AB
# AF12 becomes a (present-but-blank) "Bestamningsmetod" cell, mirroring the
# one that used to sit on row 13.
$ws.Range("AF12").Value = ""

$ws.Range("AW12").Value = $row13_AW
$ws.Range("AX12").Value = $row13_AX

# --- New row 13 (takes former row 12 content) ---
$ws.Range("A13").Value = $row12_A
$ws.Range("B13").Value = 79243
$ws.Range("Q13").Value = $row12_Q
$ws.Range("R13").Value = $row12_R
$ws.Range("S13").Value = $row12_S

$ws.Range("Z13").Clear()
$ws.Range("AB13").Clear()
$ws.Range("AF13").Clear()

$ws.Range("AW13").Value = $row12_AW
$ws.Range("AX13").Value = $row12_AX
